$wb = $excel.ActiveWorkbook
$wsTools = $wb.Worksheets.Item("Tools")
$wsSources = $wb.Worksheets.Item("Sources")

# --- Tools sheet: add a new row (10) for the Nomis Labour Market Profile tool ---
$wsTools.Range("A10").Value = "ONS Labour market profiles"
$wsTools.Range("B10").Value = "Population, employment, and qualification data by Local Authority."
$wsTools.Range("C10").Value = "<a href='https://www.nomisweb.co.uk/reports/lmp/la/contents.aspx'>ONS</a>"
$wsTools.Range("D10").Value = "Publicly available "

# Reuse the existing "left/top/wrap" formatting (same style used for the other
# source-name cells) rather than letting Excel invent a brand-new style.
$wsSources.Range("A3").Copy()
$wsTools.Range("A10").PasteSpecial(-4122)

# Row 10 wraps onto three lines, same as the other multi-line rows on the sheet.
$wsTools.Rows.Item(10).RowHeight = 43.5

# --- Sources sheet: add the Nomis Labour Market Profile link in column C of row 11 ---
$wsSources.Range("C11").Value = "Labour Market Profile - Nomis - Official Census and Labour Market Statistics (nomisweb.co.uk)"
$wsSources.Hyperlinks.Add($wsSources.Range("C11"), "https://www.nomisweb.co.uk/reports/lmp/la/contents.aspx")

# --- Restore the cursor/selection positions left behind by the edit session ---
$wsTools.Activate() | Out-Null
$wsTools.Range("C11").Select() | Out-Null

$wsSources.Activate() | Out-Null
$wsSources.Rows.Item(11).Select() | Out-Null

$wsTools.Activate() | Out-Null

Write-Output "done"
